$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 3676694.8
$ws.Range("I9").Value = 4902093
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 4902093
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = -4901924
$ws.Range("N9").Value = -838

# Row 135
$ws.Range("H135").Value = 864.1667
$ws.Range("I135").Value = 782.5
$ws.Range("J135").Value = 1517.5
$ws.Range("K135").Value = 7042.5
$ws.Range("L135").Value = 13657.5
$ws.Range("M135").Value = -4507.5
$ws.Range("N135").Value = -18727.5

# Row 138
$ws.Range("H138").Value = 3314.1143
$ws.Range("I138").Value = 3958.3635
$ws.Range("J138").Value = 3018.8333
$ws.Range("K138").Value = 11875.0905
$ws.Range("L138").Value = 9056.499899999999
$ws.Range("M138").Value = -6735.0905
$ws.Range("N138").Value = -19336.4999

$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 7493.3335
$ws.Range("I16").Value = 7493.3335
$ws.Range("K16").Value = 7493.3335
$ws.Range("M16").Value = -7206.3335

# Row 19
$ws.Range("H19").Value = 7000
$ws.Range("I19").Value = 7000
$ws.Range("K19").Value = 7000
$ws.Range("M19").Value = -6771

# Row 61
$ws.Range("H61").Value = 2118.7917
$ws.Range("I61").Value = 1577.3
$ws.Range("J61").Value = 2505.5715
$ws.Range("K61").Value = 1577.3
$ws.Range("L61").Value = 2505.5715
$ws.Range("M61").Value = -1365.3
$ws.Range("N61").Value = -2929.5715

# Row 74
$ws.Range("H74").Value = 1701.55
$ws.Range("I74").Value = 1335.5333
$ws.Range("K74").Value = 1335.5333
$ws.Range("M74").Value = -461.5333000000001

# Row 77
$ws.Range("H77").Value = 1701.55
$ws.Range("I77").Value = 1335.5333
$ws.Range("K77").Value = 6677.6665
$ws.Range("M77").Value = -2309.6665

# Row 110
$ws.Range("H110").Value = 1681.8
$ws.Range("I110").Value = 1975.5714
$ws.Range("K110").Value = 1975.5714
$ws.Range("M110").Value = 69.42859999999996

# Row 132
$ws.Range("H132").Value = 3599.4285
$ws.Range("I132").Value = 3913.1428
$ws.Range("K132").Value = 11739.4284
$ws.Range("M132").Value = -9209.428400000001

# Row 136
$ws.Range("H136").Value = 2118.7917
$ws.Range("I136").Value = 1577.3
$ws.Range("J136").Value = 2505.5715
$ws.Range("K136").Value = 4731.9
$ws.Range("L136").Value = 7516.7145
$ws.Range("M136").Value = -2181.9
$ws.Range("N136").Value = -12616.7145

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 7491.5864
$ws.Range("I134").Value = 2206.238
$ws.Range("J134").Value = 21365.625
$ws.Range("K134").Value = 6618.714
$ws.Range("L134").Value = 64096.875
$ws.Range("M134").Value = -4083.714
$ws.Range("N134").Value = -69166.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2318527.2
$ws.Range("I31").Value = 3748.3157
$ws.Range("K31").Value = 3748.3157
$ws.Range("M31").Value = -3453.3157

# Row 34
$ws.Range("H34").Value = 2318527.2
$ws.Range("I34").Value = 3748.3157
$ws.Range("K34").Value = 3748.3157
$ws.Range("M34").Value = -3546.3157

$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 629.7857
$ws.Range("J98").Value = 675.44446
$ws.Range("L98").Value = 2026.33338
$ws.Range("N98").Value = -5022.33338

# Row 121
$ws.Range("H121").Value = 14740285
$ws.Range("J121").Value = 295249.5
$ws.Range("L121").Value = 885748.5
$ws.Range("N121").Value = -888368.5

# Row 136
$ws.Range("H136").Value = 1629.6666
$ws.Range("I136").Value = 1629.6666
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4888.9998
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 211.0002000000004
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 1644.1818
$ws.Range("I22").Value = 461.6
$ws.Range("J22").Value = 2629.6667
$ws.Range("K22").Value = 461.6
$ws.Range("L22").Value = 2629.6667
$ws.Range("M22").Value = 67.39999999999998
$ws.Range("N22").Value = -3687.6667

# Row 70
$ws.Range("H70").Value = 116044.11
$ws.Range("I70").Value = 158069.08
$ws.Range("J70").Value = 6779.2
$ws.Range("K70").Value = 158069.08
$ws.Range("L70").Value = 6779.2
$ws.Range("M70").Value = -157799.08
$ws.Range("N70").Value = -7319.2

# Row 73
$ws.Range("H73").Value = 116044.11
$ws.Range("I73").Value = 158069.08
$ws.Range("J73").Value = 6779.2
$ws.Range("K73").Value = 158069.08
$ws.Range("L73").Value = 6779.2
$ws.Range("M73").Value = -157133.08
$ws.Range("N73").Value = -8651.200000000001

# Row 130
$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040

# Row 132
$ws.Range("H132").Value = 2368.7778
$ws.Range("I132").Value = 2126
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6378
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3848
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7334
$ws.Range("I7").Value = 5727.6665
$ws.Range("J7").Value = 9743.5
$ws.Range("K7").Value = 5727.6665
$ws.Range("L7").Value = 9743.5
$ws.Range("M7").Value = -5615.6665
$ws.Range("N7").Value = -9967.5

# Row 22
$ws.Range("H22").Value = 3915.4443
$ws.Range("I22").Value = 3154.875
$ws.Range("K22").Value = 3154.875
$ws.Range("M22").Value = -2859.875

# Row 27
$ws.Range("H27").Value = 3915.4443
$ws.Range("I27").Value = 3154.875
$ws.Range("K27").Value = 3154.875
$ws.Range("M27").Value = -3047.875

# Row 126
$ws.Range("H126").Value = 7334
$ws.Range("I126").Value = 5727.6665
$ws.Range("J126").Value = 9743.5
$ws.Range("K126").Value = 17182.9995
$ws.Range("L126").Value = 29230.5
$ws.Range("M126").Value = -14712.9995
$ws.Range("N126").Value = -34170.5

# Row 132
$ws.Range("H132").Value = 3944.9614
$ws.Range("I132").Value = 3372.7058
$ws.Range("K132").Value = 10118.1174
$ws.Range("M132").Value = -7588.117400000001

# Row 136
$ws.Range("H136").Value = 2245.6875
$ws.Range("J136").Value = 3426
$ws.Range("L136").Value = 10278
$ws.Range("N136").Value = -15378

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 11364716
$ws.Range("I122").Value = 1100.3684
$ws.Range("K122").Value = 3301.1052
$ws.Range("M122").Value = -851.1052

# Row 126
$ws.Range("H126").Value = 1022.6667
$ws.Range("I126").Value = 1039.5
$ws.Range("K126").Value = 3118.5
$ws.Range("M126").Value = -648.5

# Row 132
$ws.Range("H132").Value = 3740.5454
$ws.Range("I132").Value = 3532.3794
$ws.Range("J132").Value = 5249.75
$ws.Range("K132").Value = 10597.1382
$ws.Range("L132").Value = 15749.25
$ws.Range("M132").Value = -8067.138199999999
$ws.Range("N132").Value = -20809.25

# Row 136
$ws.Range("H136").Value = 205964.44
$ws.Range("I136").Value = 3628.3877
$ws.Range("K136").Value = 10885.1631
$ws.Range("M136").Value = -8335.163100000002

# Row 140
$ws.Range("H140").Value = 48447
$ws.Range("J140").Value = 48447
$ws.Range("L140").Value = 48447
$ws.Range("N140").Value = -58807
